$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New player/position/team data (rows 2-19), replacing the previous table contents.
$data = @(
    @("Austin Reaves",      "PG,SG", "Los Angeles Lakers"),
    @("Stephen Curry",      "PG,SG", "Golden State Warriors"),
    @("Darius Garland",     "PG",    "Cleveland Cavaliers"),
    @("Jaden Hardy",        "PG,SG", "Dallas Mavericks"),
    @("Keegan Murray",      "SF,PF", "Sacramento Kings"),
    @("OG Anunoby",         "SF,PF", "New York Knicks"),
    @("Karl-Anthony Towns", "PF,C",  "New York Knicks"),
    @("Kevin Durant",       "SF,PF", "Phoenix Suns"),
    @("Jarrett Allen",      "C",     "Cleveland Cavaliers"),
    @("Jalen Duren",        "C",     "Detroit Pistons"),
    @("Mark Williams",      "C",     "Charlotte Hornets"),
    @("Trey Murphy III",    "SF,PF", "New Orleans Pelicans"),
    @("Tyrese Maxey",       "PG,SG", "Philadelphia 76ers"),
    @("Dereck Lively II",   "C",     "Dallas Mavericks"),
    @("Daniel Gafford",     "PF,C",  "Dallas Mavericks"),
    @("Franz Wagner",       "SF,PF", "Orlando Magic"),
    @("Jalen Johnson",      "SF,PF", "Atlanta Hawks"),
    @("Tyrese Haliburton",  "PG,SG", "Indiana Pacers")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $row++
}
